# Update cryptocurrency price/volume figures per the refreshed data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.180.23"
$ws.Range("E2").Value = "  +1.92%  "
$ws.Range("D3").Value = "3.180.64"
$ws.Range("E3").Value = "  +3.96%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'580.09"
$ws.Range("E5").Value = "  +3.61%  "
$ws.Range("D6").Value = "'151.42"
$ws.Range("E6").Value = "  +5.54%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.178.75"
$ws.Range("E8").Value = "  +3.94%  "
$ws.Range("E9").Value = "  +3.81%  "
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "  +5.71%  "
$ws.Range("D11").Value = "'6.24"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "'0.506"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("E13").Value = "  +17.68%  "
$ws.Range("D14").Value = "'38.07"
$ws.Range("E14").Value = "  +6.12%  "
$ws.Range("D15").Value = "3.700.32"
$ws.Range("E15").Value = "  +3.89%  "
$ws.Range("D16").Value = "65.231.60"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").Value = "3.177.58"
$ws.Range("E17").Value = "  +3.62%  "
$ws.Range("D18").Value = "'7.19"
$ws.Range("E18").Value = "  +5.62%  "
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").Value = "'514.14"
$ws.Range("E20").Value = "  +7.64%  "
$ws.Range("D21").Value = "'14.91"
$ws.Range("E21").Value = "  +6.36%  "
$ws.Range("D22").Value = "'0.733"
$ws.Range("E22").Value = "  +7.02%  "
$ws.Range("D23").Value = "'15.29"
$ws.Range("E23").Value = "  +5.11%  "
$ws.Range("D24").Value = "'7.84"
$ws.Range("E24").Value = "  +3.41%  "
$ws.Range("D25").Value = "'85.49"
$ws.Range("E25").Value = "  +3.62%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "'9.03"
$ws.Range("E27").Value = "  +11.26%  "
$ws.Range("E28").Value = "  +4.61%  "
$ws.Range("E29").Value = "  +7.32%  "
$ws.Range("D30").Value = "'28.01"
$ws.Range("E30").Value = "  +6.37%  "
$ws.Range("E31").Value = "  +13.05%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'1.21"
$ws.Range("E33").Value = "  +5.60%  "
$ws.Range("D34").Value = "'6.34"
$ws.Range("E34").Value = "  +9.69%  "
$ws.Range("E35").Value = "  +6.22%  "
$ws.Range("D36").Value = "'55.71"
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("D37").Value = "'0.0905"
$ws.Range("E37").Value = "  +11.02%  "
$ws.Range("D38").Value = "'477.94"
$ws.Range("E38").Value = "  +5.78%  "
$ws.Range("E39").Value = "  +11.04%  "
$ws.Range("D40").Value = "'0.0423"
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("D41").Value = "'8.68"
$ws.Range("E41").Value = "  +4.58%  "
$ws.Range("D42").Value = "3.072.85"
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("E43").Value = "  +2.99%  "
$ws.Range("D44").Value = "'0.287"
$ws.Range("E44").Value = "  +6.90%  "
$ws.Range("D45").Value = "'2.42"
$ws.Range("E45").Value = "  +7.45%  "
$ws.Range("D46").Value = "'29.20"
$ws.Range("E46").Value = "  +3.94%  "
$ws.Range("D47").Value = "0.0₃0617"
$ws.Range("E47").Value = "  +19.55%  "
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("E50").Value = "  +7.62%  "
$ws.Range("D51").Value = "'120.63"
$ws.Range("E51").Value = "  +1.54%  "
